$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.337.76'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.686.27'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '681.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.32'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.686.19'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.29%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("E10").Value = '  -8.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.33'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.45%  '
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.63'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.310.14'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.691.15'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.423.01'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.113'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.62'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '481.52'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.91'
$ws.Range("D22").ClearFormats()
$ws.Range("E23").Value = '  -7.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.24'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.833.34'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("E26").Value = '  -9.30%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.50'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.53'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -7.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.84'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -9.33%  '
$ws.Range("E31").Value = '  -10.94%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.71%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.86'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.70%  '
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.08'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.657.92'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.47'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.34'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0936'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.26'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.24%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.954'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '161.64'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.40'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.13'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -13.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000287'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -8.91%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.11'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.57%  '
